$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 17 and 18 (existing rows 17-21 shift down to 19-23).
$ws.Rows("17:18").Insert()

# Clone formatting/merges/row-structure from the row that now holds the old row-17
# data (row 19, "VOLTAREN...") into the two new rows so every cell reuses the
# existing style indices (s6/s7/s8/s9) and merge layout instead of creating new ones.
$ws.Range("A19:N19").Copy($ws.Range("A17"))
$ws.Range("A19:N19").Copy($ws.Range("A18"))

# Restore the row heights that match the new rows' final content.
$ws.Rows(17).RowHeight = 25.5
$ws.Rows(18).RowHeight = 25.5

# --- New row 17: OTRIVIN 0.1% ADULT NASAL DROPS 15 ML ---
$ws.Range("A17").Value = 14
$ws.Range("B17").Value = "OTRIVIN 0.1% ADULT NASAL DROPS 15 ML"
$ws.Range("H17").Value = "4:0"
$ws.Range("L17").Value = 24
$ws.Range("N17").Value = "1:0"

# --- New row 18: PROXIMOL 0.4MG 40 TAB ---
$ws.Range("A18").Value = 15
$ws.Range("B18").Value = "PROXIMOL 0.4MG 40 TAB"
$ws.Range("H18").Value = "0:1"
$ws.Range("L18").Value = 34
$ws.Range("N18").Value = "1:0"

# --- Row 19 (was row 17, VOLTAREN): renumber only ---
$ws.Range("A19").Value = 16

# --- Row 20 (was row 18, WATER FOR INJECTION): renumber only ---
$ws.Range("A20").Value = 17

# --- Row 21 (was row 19, سرنجات 3 سم): renumber only ---
$ws.Range("A21").Value = 18

# --- Row 22 (was row 20, totals row): update total ---
$ws.Range("K22").Value = 1422.9500000000001

Write-Output "done"
